$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 900
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 900
$ws.Range("N9").Value = -1238
$ws.Range("M9").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 9941.727999999999
$ws.Range("I111").Value = 9078
$ws.Range("J111").Value = 12245
$ws.Range("K111").Value = 27234
$ws.Range("L111").Value = 36735
$ws.Range("M111").Value = -24167
$ws.Range("N111").Value = -42869

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 13546690
$ws.Range("I116").Value = 16931836
$ws.Range("K116").Value = 16931836
$ws.Range("M116").Value = -16928394

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7560.677
$ws.Range("I132").Value = 3300.5
$ws.Range("K132").Value = 9901.5
$ws.Range("M132").Value = -7371.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 23812414
$ws.Range("I137").Value = 2757.1428
$ws.Range("J137").Value = 47622070
$ws.Range("K137").Value = 8271.428400000001
$ws.Range("L137").Value = 142866210
$ws.Range("M137").Value = -5721.428400000001
$ws.Range("N137").Value = -142871310

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5575.985
$ws.Range("J138").Value = 6247.722
$ws.Range("L138").Value = 18743.166
$ws.Range("N138").Value = -29023.166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10287996
$ws.Range("I2").Value = 971116
$ws.Range("J2").Value = 31250976
$ws.Range("K2").Value = 971116
$ws.Range("L2").Value = 31250976
$ws.Range("M2").Value = -971003
$ws.Range("N2").Value = -31251202

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3398.0476
$ws.Range("I45").Value = 3024.6667
$ws.Range("K45").Value = 3024.6667
$ws.Range("M45").Value = -2647.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14708471
$ws.Range("I74").Value = 35715644
$ws.Range("J74").Value = 3449
$ws.Range("K74").Value = 35715644
$ws.Range("L74").Value = 3449
$ws.Range("M74").Value = -35714770
$ws.Range("N74").Value = -5197

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 14708471
$ws.Range("I77").Value = 35715644
$ws.Range("J77").Value = 3449
$ws.Range("K77").Value = 178578220
$ws.Range("L77").Value = 17245
$ws.Range("M77").Value = -178573852
$ws.Range("N77").Value = -25981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 528048.9
$ws.Range("I102").Value = 686076.0600000001
$ws.Range("K102").Value = 686076.0600000001
$ws.Range("M102").Value = -684454.0600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 10287996
$ws.Range("I116").Value = 971116
$ws.Range("J116").Value = 31250976
$ws.Range("K116").Value = 971116
$ws.Range("L116").Value = 31250976
$ws.Range("M116").Value = -968822
$ws.Range("N116").Value = -31255564

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10287996
$ws.Range("I3").Value = 971116
$ws.Range("J3").Value = 31250976
$ws.Range("K3").Value = 971116
$ws.Range("L3").Value = 31250976
$ws.Range("M3").Value = -971002
$ws.Range("N3").Value = -31251204

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 930.3333
$ws.Range("I22").Value = 914.6667
$ws.Range("K22").Value = 914.6667
$ws.Range("M22").Value = -741.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 761714.75
$ws.Range("I94").Value = 761714.75
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 761714.75
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -761263.75
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 113427420
$ws.Range("I99").Value = 113427420
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 113427420
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -113425922
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 57694270
$ws.Range("I105").Value = 57694270
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 57694270
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -57692523
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3205.5
$ws.Range("I107").Value = 3205.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3205.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1285.5
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2293.3447
$ws.Range("I134").Value = 2200
$ws.Range("K134").Value = 6600
$ws.Range("M134").Value = -4065

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45457616
$ws.Range("I31").Value = 52634012
$ws.Range("K31").Value = 52634012
$ws.Range("M31").Value = -52633717

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 45457616
$ws.Range("I34").Value = 52634012
$ws.Range("K34").Value = 52634012
$ws.Range("M34").Value = -52633810

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 49833.668
$ws.Range("J62").Value = 72499.664
$ws.Range("L62").Value = 72499.664
$ws.Range("N62").Value = -73747.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 49833.668
$ws.Range("J65").Value = 72499.664
$ws.Range("L65").Value = 362498.32
$ws.Range("N65").Value = -368738.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 32243.75
$ws.Range("J97").Value = 32243.75
$ws.Range("L97").Value = 32243.75
$ws.Range("N97").Value = -34225.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1818927.9
$ws.Range("I107").Value = 3030897
$ws.Range("K107").Value = 3030897
$ws.Range("M107").Value = -3028977

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2211583.2
$ws.Range("I4").Value = 480939.25
$ws.Range("K4").Value = 1442817.75
$ws.Range("M4").Value = -1442705.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1342.8667
$ws.Range("I5").Value = 801.1667
$ws.Range("J5").Value = 1704
$ws.Range("K5").Value = 2403.5001
$ws.Range("L5").Value = 5112
$ws.Range("M5").Value = -2291.5001
$ws.Range("N5").Value = -5336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 96.333336
$ws.Range("I8").Value = 96.333336
$ws.Range("K8").Value = 289.000008
$ws.Range("M8").Value = -150.000008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 250057470
$ws.Range("J37").Value = 250057470
$ws.Range("L37").Value = 750172410
$ws.Range("N37").Value = -750172634

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1342.8667
$ws.Range("I135").Value = 801.1667
$ws.Range("J135").Value = 1704
$ws.Range("K135").Value = 7210.5003
$ws.Range("L135").Value = 15336
$ws.Range("M135").Value = -4675.5003
$ws.Range("N135").Value = -20406

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 32695694
$ws.Range("I137").Value = 125001170
$ws.Range("J137").Value = 5004049.5
$ws.Range("K137").Value = 375003510
$ws.Range("L137").Value = 15012148.5
$ws.Range("M137").Value = -374998410
$ws.Range("N137").Value = -15022348.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4502
$ws.Range("I5").Value = 4502
$ws.Range("K5").Value = 4502
$ws.Range("M5").Value = -4390

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 11999
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 6803454.5
$ws.Range("I107").Value = 9524156
$ws.Range("J107").Value = 1699.5
$ws.Range("K107").Value = 9524156
$ws.Range("L107").Value = 1699.5
$ws.Range("M107").Value = -9522236
$ws.Range("N107").Value = -5539.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3452.4
$ws.Range("J113").Value = 3739.8
$ws.Range("L113").Value = 3739.8
$ws.Range("N113").Value = -8079.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 531422.9399999999
$ws.Range("I122").Value = 1576284.1
$ws.Range("K122").Value = 4728852.300000001
$ws.Range("M122").Value = -4726402.300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 22227992
$ws.Range("J40").Value = 33339634
$ws.Range("L40").Value = 33339634
$ws.Range("N40").Value = -33339906

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1796.9231
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 557969.5
$ws.Range("I100").Value = 751382.5600000001
$ws.Range("K100").Value = 1502765.12
$ws.Range("M100").Value = -1502224.12

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1049.579
$ws.Range("J113").Value = 1258.7778
$ws.Range("L113").Value = 3776.3334
$ws.Range("N113").Value = -8116.3334
